$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) HEADER sheet: reorder the key/value rows and add a DOMAIN row.
#    before: A3 SOURCE_ORG/-        A4 SOURCE_PERSON/-
#            A5 CATEGORY/IFDAT      A6 SUB_CATEGORY/OTHER_KEY
#    after:  A3 DOMAIN/IFDAT        A4 CATEGORY/OTHER_KEY
#            A5 SOURCE_ORG/-        A6 SOURCE_PERSON/-
# ---------------------------------------------------------------------------
$header = $wb.Worksheets.Item("HEADER")
$header.Range("A3").Value = "DOMAIN"
$header.Range("B3").Value = "IFDAT"
$header.Range("A4").Value = "CATEGORY"
$header.Range("B4").Value = "OTHER_KEY"
$header.Range("A5").Value = "SOURCE_ORG"
$header.Range("B5").Value = ""
$header.Range("A6").Value = "SOURCE_PERSON"
$header.Range("B6").Value = ""

# ---------------------------------------------------------------------------
# 2) DIVIDEND_ sheet: insert a new "-" entry at the top of the TYP code list
#    (col A), pushing the existing 14 codes down one row (col B/C untouched).
#    Use Range.Copy so the cells stay text-typed (no NumberFormat/style churn).
# ---------------------------------------------------------------------------
$dividendCodes = $wb.Worksheets.Item("DIVIDEND_")
$dividendCodes.Range("A1:A14").Copy($dividendCodes.Range("A2:A15"))
$dividendCodes.Range("A1").Value = "-"

# ---------------------------------------------------------------------------
# 3) DIVIDEND sheet: the TYP column's validation list now covers the extra
#    row (A1:A15 instead of A1:A14).
# ---------------------------------------------------------------------------
$dividend = $wb.Worksheets.Item("DIVIDEND")
$dividend.Range("C4:C20").Validation.Formula1 = "='DIVIDEND_'!`$A`$1:`$A`$15"
